$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# First, move the "Bring Up" / "Test " rows down from 25/26 to 30/31
# (read values via the Value() method call form, which this runtime requires
# for reliable reads; the Value property getter returns a descriptor string)
$bringUp = $ws.Range("A25").Value()
$testVal = $ws.Range("A26").Value()

$ws.Range("A25").ClearContents()
$ws.Range("A26").ClearContents()

$ws.Range("A30").Value = $bringUp
$ws.Range("A31").Value = $testVal

# Add Number (F column) values for rows 17, 18, 20, 21, 22, 23
$ws.Range("F17").Value = 3
$ws.Range("F18").Value = 8
$ws.Range("F20").Value = 4
$ws.Range("F21").Value = 3
$ws.Range("F22").Value = 2
$ws.Range("F23").Value = 4

# New rows 25-28: Eagle related sub tasks under "Realization"
$ws.Range("B25").Value = "Eagle library Resarch"
$ws.Range("F25").Value = 4

$ws.Range("B26").Value = "Eagle library build"
$ws.Range("F26").Value = 3

$ws.Range("B27").Value = "Eagle schematic"
$ws.Range("F27").Value = 3

$ws.Range("B28").Value = "Eagle board"
$ws.Range("F28").Value = 2

# Update selection to match the new active cell
$ws.Range("G28").Select()
